$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(2,  -0.7929665312411044, 0.4332962080770135),
    @(3,  -0.3291314904302632, 0.7440752289085624),
    @(4,  -1.48222548100006,   0.1474903822841229),
    @(5,  -0.3177875479579046, 0.7525896968259169),
    @(6,   0.4141323824909867, 0.6813782098136179),
    @(7,  -0.4132982357630789, 0.6819833074473352),
    @(8,   0.2364224499526925, 0.8145234463099507),
    @(9,  -0.9369859605422933, 0.3553741610524392),
    @(10, -0.07576364138646027, 0.9400511985209494),
    @(11,  1.025874795791025,  0.3121942834871412)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
}
